$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YES PHN")

$ws.Range("AG2").Value = 99
$ws.Range("AH2").Value = 99
$ws.Range("AG3").Value = 99
$ws.Range("AH3").Value = 99
